$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "创建日期" (creation date) column for rows 23-26 from 2025/6/14 to 2025/6/15
$ws.Range("F23:F26").Value = "2025/6/15"

# Reflect the user's selection/scroll position at the time of the edit
$ws.Range("F23:F26").Select()
$ws.Application.ActiveWindow.ScrollRow = 23
